$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.809.81"
$ws.Range("E2").Value = "  +0.24%  "

$ws.Range("D3").Value = "2.410.12"
$ws.Range("E3").Value = "  -0.21%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'551.06"
$ws.Range("E5").Value = "  -0.54%  "

$ws.Range("D6").Value = "'136.97"
$ws.Range("E6").Value = "  -0.83%  "

$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").Value = "'0.591"
$ws.Range("E8").Value = "  +3.78%  "

$ws.Range("E9").Value = "  -1.59%  "

$ws.Range("D10").Value = "'5.67"
$ws.Range("E10").Value = "  -2.75%  "

$ws.Range("D11").Value = "'0.148"
$ws.Range("E11").Value = "  -0.95%  "

$ws.Range("E12").Value = "  -1.72%  "

$ws.Range("D13").Value = "'25.24"
$ws.Range("E13").Value = "  +1.93%  "

$ws.Range("D14").Value = "2.838.12"
$ws.Range("E14").Value = "  -0.38%  "

$ws.Range("D15").Value = "59.749.28"
$ws.Range("E15").Value = "  +0.28%  "

$ws.Range("D17").Value = "2.420.87"
$ws.Range("E17").Value = "  -0.09%  "

$ws.Range("D18").Value = "'11.31"
$ws.Range("E18").Value = "  -0.71%  "

$ws.Range("E19").Value = "  -0.55%  "

$ws.Range("D20").Value = "'328.45"
$ws.Range("E20").Value = "  -1.49%  "

$ws.Range("D21").Value = "'6.67"
$ws.Range("E21").Value = "  -3.53%  "

$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.01%  "

$ws.Range("D23").Value = "'66.07"
$ws.Range("E23").Value = "  +2.44%  "

$ws.Range("E24").Value = "  +2.39%  "

$ws.Range("D25").Value = "'8.61"
$ws.Range("E25").Value = "  +0.47%  "

$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("D27").Value = "'1.37"
$ws.Range("E27").Value = "  +0.55%  "

$ws.Range("D28").Value = "0.0₃0769"
$ws.Range("E28").Value = "  -1.93%  "

$ws.Range("E29").Value = "  -2.35%  "

$ws.Range("D30").Value = "'169.26"
$ws.Range("E30").Value = "  -0.83%  "

$ws.Range("D31").Value = "'6.02"
$ws.Range("E31").Value = "  -3.99%  "

$ws.Range("E32").Value = "  -0.38%  "

$ws.Range("E33").Value = "  -1.70%  "

$ws.Range("E35").Value = "  -0.49%  "

$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("D37").Value = "'4.17"
$ws.Range("E37").Value = "  -1.43%  "

$ws.Range("E38").Value = "  -2.08%  "

$ws.Range("D39").Value = "'321.03"
$ws.Range("E39").Value = "  +3.00%  "

$ws.Range("D40").Value = "'0.403"
$ws.Range("E40").Value = "  -4.59%  "

$ws.Range("D41").Value = "'3.64"
$ws.Range("E41").Value = "  -2.39%  "

$ws.Range("D42").Value = "'139.83"
$ws.Range("E42").Value = "  -2.30%  "

$ws.Range("D43").Value = "'0.0966"
$ws.Range("E43").Value = "  +0.22%  "

$ws.Range("D44").Value = "'19.48"
$ws.Range("E44").Value = "  +1.72%  "

$ws.Range("E45").Value = "  -1.82%  "

$ws.Range("D46").Value = "'0.577"
$ws.Range("E46").Value = "  +1.36%  "

$ws.Range("D47").Value = "'0.402"
$ws.Range("E47").Value = "  -1.86%  "

$ws.Range("D48").Value = "'0.0223"
$ws.Range("E48").Value = "  -1.43%  "

$ws.Range("D49").Value = "'11.03"
$ws.Range("E49").Value = "  -0.16%  "

$ws.Range("E50").Value = "  -3.52%  "

$ws.Range("D51").Value = "'4.68"
$ws.Range("E51").Value = "  -0.98%  "
